# Add the new Covid cumulative-deaths rows (76-83, dates 2020-12-29 .. 2021-01-06)
# to Sheet1, matching the style already used for the "Date" column (column A),
# and move the active selection to C69 as in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# ---------------------------------------------------------------------------
# 1. Bake the new cell style (date number format + left-aligned) on a scratch
#    cell that is well outside the used range, so the style gets created
#    cleanly (column A already carries a column-level date style, and mutating
#    alignment directly on a cell that inherits that style causes the engine
#    to emit a redundant custom number-format entry). We then copy the style
#    onto the new date cells and clean the scratch cell back up.
# ---------------------------------------------------------------------------
$scratch = $ws.Cells.Item(1, 50)
$scratch.HorizontalAlignment = -4131   # xlLeft
$scratch.NumberFormat = "mm-dd-yy"     # maps to built-in date format (numFmtId 14)

$scratch.Copy()
$ws.Range("A76:A83").PasteSpecial(-4122)   # xlPasteFormats

$scratch.Clear()
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Write the new rows of data.
# ---------------------------------------------------------------------------
$newRows = @(
    @(76, 44194, 2065, 507, 2572),
    @(77, 44195, 2138, 531, 2669),
    @(78, 44196, 2250, 569, 2819),
    @(79, 44198, 2317, 591, 2908),
    @(80, 44199, 2521, 610, 3131),
    @(81, 44200, 2603, 622, 3225),
    @(82, 44201, 2657, 627, 3284),
    @(83, 44202, 2717, 632, 3349)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# ---------------------------------------------------------------------------
# 3. Update view state: scroll position and current selection.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 41
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C69").Select()
